$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.771.80'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '2.630.35'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.11'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("E6").Value = '  +1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.68'
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").Value = '3.105.12'
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '63.705.10'
$ws.Range("D17").Value = '2.657.46'
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.17'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.71'
$ws.Range("E19").Value = '  +2.50%  '
$ws.Range("E20").Value = '  -2.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.47'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.32'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("E24").Value = '  +9.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000112'
$ws.Range("E25").Value = '  +2.91%  '
$ws.Range("E26").Value = '  +3.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.22'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '581.25'
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.28'
$ws.Range("E29").Value = '  +4.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.160'
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.65'
$ws.Range("E34").Value = '  +2.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.47'
$ws.Range("E35").Value = '  +2.89%  '
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.78'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("E39").Value = '  +2.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.89'
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E41").Value = '  +8.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '162.92'
$ws.Range("E43").Value = '  +3.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.07'
$ws.Range("E44").Value = '  +5.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.92'
$ws.Range("E45").Value = '  -0.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0587'
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.634'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0248'
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("D50").Value = '0.0₆0238'
$ws.Range("E50").Value = '  +2.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.788'
$ws.Range("E51").Value = '  +1.86%  '
